$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.964.49"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "3.133.95"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.81"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.83"
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.130.76"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +4.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.30"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.415"
$ws.Range("E12").Value = "  +4.30%  "
$ws.Range("D13").Value = "3.670.71"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.137"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.61"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "58.092.12"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "3.138.54"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.04"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.73"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.14"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "360.06"
$ws.Range("E22").Value = "  +3.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.04"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.506"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("E27").Value = "  +2.57%  "
$ws.Range("D28").Value = "0.0₃0874"
$ws.Range("E28").Value = "  -4.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.32"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.44"
$ws.Range("E32").Value = "  +1.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.10"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.89"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.08"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.27"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("E39").Value = "  +2.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0672"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("D41").Value = "2.502.50"
$ws.Range("E41").Value = "  +6.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "37.73"
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").Value = "3.175.49"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0267"
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.992"
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.08"
$ws.Range("E49").Value = "  +0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.82"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.740"
$ws.Range("E51").Value = "  -3.17%  "

# Row 42 and 43 swap: Mantle <-> Filecoin
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.01"
$ws.Range("E42").Value = "  -4.38%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.699"
$ws.Range("E43").Value = "  -0.55%  "
